# ----------------------------------------------------------------------
# Replace old Constructor/Encapsulation questions (rows 21-23) with new
# Exceptions / Thread questions, per commit:
# "Few questions added for Exception and Threads in MTT for all the Batches"
# ----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlHAlignLeft   = -4131
$xlHAlignRight  = -4152
$xlHAlignCenter = -4108
$xlHAlignGeneral = 1

# ---- Row 21 ----
$t1 = @'
Exceptions
'@
$ws.Range("A21").Value2 = $t1
$ws.Range("A21").WrapText = $false
$ws.Range("A21").HorizontalAlignment = $xlHAlignGeneral

$t2 = @'
Code
'@
$ws.Range("B21").Value2 = $t2
$ws.Range("B21").Font.Size = 11
$ws.Range("B21").WrapText = $false
$ws.Range("B21").HorizontalAlignment = $xlHAlignLeft

$t3 = @'
Single Choice
'@
$ws.Range("C21").Value2 = $t3
$ws.Range("C21").Font.Size = 11
$ws.Range("C21").WrapText = $false
$ws.Range("C21").HorizontalAlignment = $xlHAlignCenter

$t4 = @'
class Test extends Exception { }
class Main {
   public static void main(String args[]) { 
      try {
         throw new Test();
      }
      catch(Test t) {
         System.out.println("Got the Test Exception");
      }
      finally {
         System.out.println("Inside finally block ");
      }
  }
}
'@
$ws.Range("D21").Value2 = $t4
$ws.Range("D21").WrapText = $true
$ws.Range("D21").HorizontalAlignment = $xlHAlignGeneral

$t5 = @'
Got the Test Exception
Inside finally block 
'@
$ws.Range("E21").Value2 = $t5
$ws.Range("E21").WrapText = $true
$ws.Range("E21").HorizontalAlignment = $xlHAlignGeneral

$t6 = @'
Got the Test Exception
'@
$ws.Range("F21").Value2 = $t6
$ws.Range("F21").WrapText = $false
$ws.Range("F21").HorizontalAlignment = $xlHAlignGeneral

$t7 = @'
Inside finally block 
'@
$ws.Range("G21").Value2 = $t7
$ws.Range("G21").WrapText = $false
$ws.Range("G21").HorizontalAlignment = $xlHAlignGeneral

$t8 = @'
Compiler Error
'@
$ws.Range("H21").Value2 = $t8
$ws.Range("H21").WrapText = $false
$ws.Range("H21").HorizontalAlignment = $xlHAlignGeneral

$ws.Range("I21").Value2 = 1
$ws.Range("I21").WrapText = $false
$ws.Range("I21").HorizontalAlignment = $xlHAlignRight

$ws.Rows.Item(21).RowHeight = 255

# ---- Row 22 ----
$t9 = @'
Thread
'@
$ws.Range("A22").Value2 = $t9
$ws.Range("A22").WrapText = $false
$ws.Range("A22").HorizontalAlignment = $xlHAlignGeneral

$t10 = @'
Code
'@
$ws.Range("B22").Value2 = $t10
$ws.Range("B22").Font.Size = 11
$ws.Range("B22").WrapText = $false
$ws.Range("B22").HorizontalAlignment = $xlHAlignLeft

$t11 = @'
Single Choice
'@
$ws.Range("C22").Value2 = $t11
$ws.Range("C22").Font.Size = 11
$ws.Range("C22").WrapText = $false
$ws.Range("C22").HorizontalAlignment = $xlHAlignGeneral

$t12 = @'
class MultithreadedPrograming
    {
        public static void main(String args[])
        {
            Thread t = Thread.currentThread();
            t.setName("New Thread");
            System.out.println(t);        
        }
    }
'@
$ws.Range("D22").Value2 = $t12
$ws.Range("D22").WrapText = $true
$ws.Range("D22").HorizontalAlignment = $xlHAlignGeneral

$t13 = @'
Thread[5,main].
'@
$ws.Range("E22").Value2 = $t13
$ws.Range("E22").WrapText = $false
$ws.Range("E22").HorizontalAlignment = $xlHAlignGeneral

$t14 = @'
Thread[New Thread,5].
'@
$ws.Range("F22").Value2 = $t14
$ws.Range("F22").WrapText = $false
$ws.Range("F22").HorizontalAlignment = $xlHAlignGeneral

$t15 = @'
Thread[main,5,main].
'@
$ws.Range("G22").Value2 = $t15
$ws.Range("G22").WrapText = $false
$ws.Range("G22").HorizontalAlignment = $xlHAlignGeneral

$t16 = @'
Thread[New Thread,5,main].
'@
$ws.Range("H22").Value2 = $t16
$ws.Range("H22").WrapText = $false
$ws.Range("H22").HorizontalAlignment = $xlHAlignGeneral

$ws.Range("I22").Value2 = 4
$ws.Range("I22").WrapText = $false
$ws.Range("I22").HorizontalAlignment = $xlHAlignGeneral

$ws.Rows.Item(22).RowHeight = 135

# ---- Row 23 ----
$t17 = @'
Exceptions
'@
$ws.Range("A23").Value2 = $t17
$ws.Range("A23").WrapText = $false
$ws.Range("A23").HorizontalAlignment = $xlHAlignLeft

$t18 = @'
Code
'@
$ws.Range("B23").Value2 = $t18
$ws.Range("B23").Font.Size = 11
$ws.Range("B23").WrapText = $false
$ws.Range("B23").HorizontalAlignment = $xlHAlignLeft

$t19 = @'
Single Choice
'@
$ws.Range("C23").Value2 = $t19
$ws.Range("C23").Font.Size = 11
$ws.Range("C23").WrapText = $false
$ws.Range("C23").HorizontalAlignment = $xlHAlignLeft

$t20 = @'
class Test
{
    public static void main (String[] args)
    {
        try
        {
            int a = 0;
            System.out.println ("a = " + a);
            int b = 20 / a;
            System.out.println ("b = " + b);
        }
        catch(ArithmeticException e)
        {
            System.out.println ("Divide by zero error");
        }
        finally
        {
            System.out.println ("inside the finally block");
        }
    }
}
'@
$ws.Range("D23").Value2 = $t20
$ws.Range("D23").WrapText = $true
$ws.Range("D23").HorizontalAlignment = $xlHAlignLeft

$t21 = @'
Compile error
'@
$ws.Range("E23").Value2 = $t21
$ws.Range("E23").WrapText = $false
$ws.Range("E23").HorizontalAlignment = $xlHAlignLeft

$t22 = @'
a = 0
Divide by zero error
inside the finally block
'@
$ws.Range("F23").Value2 = $t22
$ws.Range("F23").WrapText = $true
$ws.Range("F23").HorizontalAlignment = $xlHAlignLeft

$t23 = @'
A = 0
'@
$ws.Range("G23").Value2 = $t23
$ws.Range("G23").WrapText = $false
$ws.Range("G23").HorizontalAlignment = $xlHAlignLeft

$t24 = @'
inside the finally block
'@
$ws.Range("H23").Value2 = $t24
$ws.Range("H23").WrapText = $false
$ws.Range("H23").HorizontalAlignment = $xlHAlignLeft

$ws.Range("I23").Value2 = 2
$ws.Range("I23").WrapText = $false
$ws.Range("I23").HorizontalAlignment = $xlHAlignRight

$ws.Rows.Item(23).RowHeight = 375

# ---- Update active selection / scroll position ----
$ws.Activate()
$ws.Range("A23:I23").Select()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1

